$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45944
$ws.Range("B2").Value = 0.007679953205691472
$ws.Range("C2").Value = 1.626911367747524
$ws.Range("D2").Value = 0.02282632902275583
$ws.Range("E2").Value = 0.00005898168124161073
$ws.Range("F2").Value = 29
$ws.Range("G2").Value = 3492.988802624848
$ws.Range("H2").Value = 59.10151269320311
$ws.Range("I2").Value = 49.10089788410074
$ws.Range("J2").Value = 0.8710329042982154
